$wb = $excel.ActiveWorkbook

# --- Environment_PartnsCom (2nd sheet) : insert a new row 5 with the
#     "SOI Testing Farmer" / envUsernameNameITTQA entry, pushing the
#     existing rows 5-12 down to 6-13. ---
$ws2 = $wb.Worksheets.Item(2)

$ws2.Rows.Item(5).Insert() | Out-Null
$ws2.Range("A5").Value = "envUsernameNameITTQA"
$ws2.Range("B5").Value = "SOI Testing Farmer"

# column A widened a bit to fit the new values (best-fit recalculated by Excel)
$ws2.Columns.Item(1).ColumnWidth = 23.45

# hyperlinks: collection doesn't auto-shift on row insert, so rebuild them in
# their original order (same target addresses) so the relationship ids line
# up again with the (now shifted) cells that own them.
$ws2.Hyperlinks.Delete() | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("B3"), "mailto:soi.testing.crew@gmail.com.farmer") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("B7"), "mailto:soi.testing.crew@gmail.com.hunter") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("B9"), "mailto:soi.testing.crew@gmail.com.nonexclusivehunter") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("B2"), "https://prxitt-proximus.cs127.force.com/SalesforceforPartners") | Out-Null

# move the active tab / selection to this sheet
$ws2.Activate()
$ws2.Range("E15").Select() | Out-Null
